$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-19 (row 1 is the header and is unchanged)
$data = @(
    @("Stephon Castle", "PG,SG", "San Antonio Spurs"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("Goga Bitadze", "C", "Orlando Magic"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Santi Aldama", "PF,C", "Memphis Grizzlies"),
    @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers"),
    @("Alperen Sengün", "C", "Houston Rockets"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers"),
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
